$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> Kolkata Knight Riders vs Royal Challengers Bengaluru (was row 4)
$ws.Range("A2").Value = "22-03-2025"
$ws.Range("B2").Value = "Kolkata Knight Riders vs Royal Challengers Bengaluru"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Kolkata Knight Riders"

# Row 3 -> Sunrisers Hyderabad vs Rajasthan Royals (was row 2)
$ws.Range("A3").Value = "23-03-2025"
$ws.Range("B3").Value = "Sunrisers Hyderabad vs Rajasthan Royals"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "Rajasthan Royals"

# Row 4 -> Chennai Super Kings vs Mumbai Indians (was row 3)
$ws.Range("A4").Value = "23-03-2025"
$ws.Range("B4").Value = "Chennai Super Kings vs Mumbai Indians"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "Mumbai Indians"
$ws.Range("E4").Value = "Mumbai Indians"
